# Update schedule data on main sheet (LP1912)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Última actualización: 11:54:18"
$ws.Range("A3").Value = "Total filas: 151"
$ws.Range("A23").Value = "06:17:28"
$ws.Range("C23").Value = "16_SANTA ANA"
$ws.Range("D23").Value = 64
$ws.Range("A24").Value = "05:57:13"
$ws.Range("C24").Value = "23_HERNANDEZ"
$ws.Range("D24").Value = 84
$ws.Range("C40").Value = "11_ETCHEVERRY"
$ws.Range("C41").Value = "15_ABASTO"
$ws.Range("A99").Value = "10:05:51"
$ws.Range("C99").Value = "16_SANTA ANA"
$ws.Range("D99").Value = 86
$ws.Range("A100").Value = "11:11:33"
$ws.Range("C100").Value = "215C_EL PATO"
$ws.Range("D100").Value = 20
$ws.Range("A111").Value = "10:50:41"
$ws.Range("C111").Value = "23_HERNANDEZ"
$ws.Range("D111").Value = 64
$ws.Range("A112").Value = "11:54:18"
$ws.Range("C112").Value = "15X38_ABASTO"
$ws.Range("D112").Value = 0
$ws.Range("A113").Value = "11:52:01"
$ws.Range("B113").Value = "11:54"
$ws.Range("C113").Value = "225_GOMEZ"
$ws.Range("D113").Value = 2
$ws.Range("A114").Value = "11:34:59"
$ws.Range("B114").Value = "11:57"
$ws.Range("D114").Value = 23
$ws.Range("A115").Value = "10:05:51"
$ws.Range("B115").Value = "11:58"
$ws.Range("C115").Value = "17_ROMERO"
$ws.Range("D115").Value = 113
$ws.Range("A116").Value = "10:37:52"
$ws.Range("B116").Value = "12:05"
$ws.Range("D116").Value = 88
$ws.Range("A117").Value = "11:47:17"
$ws.Range("B117").Value = "12:06"
$ws.Range("C117").Value = "11_ETCHEVERRY"
$ws.Range("D117").Value = 19
$ws.Range("C118").Value = "16_P MOR-SANTA ANA"
$ws.Range("A119").Value = "11:34:59"
$ws.Range("B119").Value = "12:09"
$ws.Range("D119").Value = 35
$ws.Range("C120").Value = "15_ABASTO"
$ws.Range("B121").Value = "12:10"
$ws.Range("C121").Value = "16_P MOR-SANTA ANA"
$ws.Range("D121").Value = 93
$ws.Range("A122").Value = "10:37:52"
$ws.Range("B122").Value = "12:16"
$ws.Range("D122").Value = 99
$ws.Range("A123").Value = "11:11:33"
$ws.Range("B123").Value = "12:17"
$ws.Range("C123").Value = "10_OLMOS"
$ws.Range("D123").Value = 66
$ws.Range("A124").Value = "10:37:52"
$ws.Range("B124").Value = "12:21"
$ws.Range("D124").Value = 104
$ws.Range("A125").Value = "11:11:33"
$ws.Range("B125").Value = "12:22"
$ws.Range("C125").Value = "215C_EL PATO"
$ws.Range("D125").Value = 71
$ws.Range("A126").Value = "11:47:17"
$ws.Range("C126").Value = "23_HERNANDEZ"
$ws.Range("D126").Value = 45
$ws.Range("A127").Value = "10:37:52"
$ws.Range("B127").Value = "12:32"
$ws.Range("D127").Value = 115
$ws.Range("A128").Value = "11:47:17"
$ws.Range("C128").Value = "14_ABASTO"
$ws.Range("D128").Value = 46
$ws.Range("A129").Value = "11:34:59"
$ws.Range("B129").Value = "12:33"
$ws.Range("D129").Value = 59
$ws.Range("A130").Value = "10:37:52"
$ws.Range("B130").Value = "12:34"
$ws.Range("C130").Value = "15_ABASTO"
$ws.Range("D130").Value = 117
$ws.Range("A131").Value = "11:11:33"
$ws.Range("C131").Value = "23_HERNANDEZ"
$ws.Range("D131").Value = 84
$ws.Range("A132").Value = "11:34:59"
$ws.Range("B132").Value = "12:35"
$ws.Range("D132").Value = 61
$ws.Range("A133").Value = "10:50:41"
$ws.Range("C133").Value = "27_EL RETIRO"
$ws.Range("D133").Value = 106
$ws.Range("A134").Value = "11:34:59"
$ws.Range("B134").Value = "12:36"
$ws.Range("C134").Value = "23_HERNANDEZ"
$ws.Range("D134").Value = 62
$ws.Range("A135").Value = "11:47:17"
$ws.Range("C135").Value = "27_EL RETIRO"
$ws.Range("D135").Value = 50
$ws.Range("A136").Value = "11:52:01"
$ws.Range("B136").Value = "12:37"
$ws.Range("C136").Value = "23_HERNANDEZ"
$ws.Range("D136").Value = 45
$ws.Range("C137").Value = "15X38_ABASTO"
$ws.Range("C138").Value = "16_SANTA ANA"
$ws.Range("A139").Value = "11:34:59"
$ws.Range("B139").Value = "12:47"
$ws.Range("D139").Value = 73
$ws.Range("A140").Value = "11:47:17"
$ws.Range("C140").Value = "14_ABASTO"
$ws.Range("D140").Value = 61
$ws.Range("A141").Value = "11:11:33"
$ws.Range("C141").Value = "15X38_ABASTO"
$ws.Range("D141").Value = 97
$ws.Range("A142").Value = "10:50:41"
$ws.Range("B142").Value = "12:48"
$ws.Range("C142").Value = "16_SANTA ANA"
$ws.Range("D142").Value = 118
$ws.Range("A143").Value = "11:11:33"
$ws.Range("B143").Value = "13:02"
$ws.Range("C143").Value = "11_ETCHEVERRY"
$ws.Range("D143").Value = 111
$ws.Range("A144").Value = "11:34:59"
$ws.Range("C144").Value = "215C_EL PATO"
$ws.Range("D144").Value = 89
$ws.Range("B145").Value = "13:03"
$ws.Range("C145").Value = "11_ETCHEVERRY"
$ws.Range("D145").Value = 76
$ws.Range("A146").Value = "11:47:17"
$ws.Range("B146").Value = "13:04"
$ws.Range("C146").Value = "215C_EL PATO"
$ws.Range("D146").Value = 77
$ws.Range("A147").Value = "11:34:59"
$ws.Range("B147").Value = "13:12"
$ws.Range("D147").Value = 98
$ws.Range("A148").Value = "11:47:17"
$ws.Range("B148").Value = "13:13"
$ws.Range("C148").Value = "16_SANTA ANA"
$ws.Range("D148").Value = 86
$ws.Range("A149").Value = "11:34:59"
$ws.Range("B149").Value = "13:16"
$ws.Range("D149").Value = 102
$ws.Range("A150").Value = "11:47:17"
$ws.Range("B150").Value = "13:17"
$ws.Range("C150").Value = "10_OLMOS"
$ws.Range("D150").Value = 90
$ws.Range("A151").Value = "11:54:18"
$ws.Range("B151").Value = "13:22"
$ws.Range("C151").Value = "23_HERNANDEZ"
$ws.Range("D151").Value = 88
$ws.Range("B152").Value = "13:24"
$ws.Range("C152").Value = "16_P MOR-SANTA ANA"
$ws.Range("D152").Value = 110
$ws.Range("B153").Value = "13:25"
$ws.Range("C153").Value = "16_P MOR-SANTA ANA"
$ws.Range("D153").Value = 98
$ws.Range("A154").Value = "11:34:59"
$ws.Range("B154").Value = "13:32"
$ws.Range("C154").Value = "215A_EL PATO"
$ws.Range("D154").Value = 118
$ws.Range("A155").Value = "11:47:17"
$ws.Range("B155").Value = "13:33"
$ws.Range("C155").Value = "215A_EL PATO"
$ws.Range("D155").Value = 106
$ws.Range("E155").Value = "LP1912"
$ws.Range("A156").Value = "11:52:01"
$ws.Range("B156").Value = "13:47"
$ws.Range("C156").Value = "225_GOMEZ"
$ws.Range("D156").Value = 115
$ws.Range("E156").Value = "LP1912"

# Update "Ultima actualizacion" timestamp on the other two sheets
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 11:54:18"
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "Última actualización: 11:54:18"
